$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("3:3").Insert()

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2023-11-27"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "케이엔에스"
$ws.Range("C3").Value = "신영"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2023-11-30"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2023-12-06"
$ws.Range("E3").ClearFormats()

$ws.Range("F3").Value = 17250000
$ws.Range("G3").Value = 750000
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = 19000
$ws.Range("J3").Value = 22000
$ws.Range("K3").Value = "-"
$ws.Range("L3").Value = 23000
$ws.Range("M3").Value = "-"
$ws.Range("N3").Value = "-"
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "-"
$ws.Range("Q3").Value = "-"
$ws.Range("R3").Value = "1450.7 : 1"
$ws.Range("S3").Value = "-"
$ws.Range("T3").Value = "-"
